$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading 1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$afterTitle = $titlePara.Range
$afterTitle.Collapse(0)               # wdCollapseEnd
$afterTitle.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaRange = $metaPara.Range.Duplicate
$metaRange.Collapse(1)                # wdCollapseStart

$metaRange.InsertAfter("Meta description")
$metaRange.Font.Bold = 1

$afterBold = $d.Range($metaRange.End, $metaRange.End)
$afterBold.InsertAfter(": Read our review of Aztec Bonanza, the high-variance slot game with cascading wins and a 19,000x jackpot payout. Play for free and win big today!")
$afterBold.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) Near the end of the document:
#      - delete the bold "Play Aztec Bonanza Free..." paragraph
#      - replace the italic paragraph's text with the new image-prompt text
# ---------------------------------------------------------------------------
$boldPara = $null
$italicPara = $null
for ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t.StartsWith("Play Aztec Bonanza Free - Exciting Jackpot Payout")) {
        $boldPara = $p
    }
    if ($t.StartsWith("Read our review of Aztec Bonanza")) {
        $italicPara = $p
    }
    if ($boldPara -ne $null -and $italicPara -ne $null) {
        break
    }
}

if ($italicPara -eq $null -or $boldPara -eq $null) {
    Write-Output ("ERROR: could not locate target paragraphs. boldPara=" + $boldPara + " italicPara=" + $italicPara)
} else {
    # Replace the italic paragraph's text (but keep its paragraph mark/formatting)
    # with the new image-generation prompt text.
    $italicTextRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End - 1)
    $italicTextRange.Text = "Create a feature image for Aztec Bonanza that features a happy Maya warrior wearing glasses in a cartoon style. The image should be lively and exciting to match the adventurous theme of the game. You can incorporate elements like gemstones, the totem pole, or the temple in the background to give the image an Aztec touch. Make sure to highlight the cascading game mechanics and the free spins and giant symbol features to entice players. The overall image should convey a feeling of excitement and adventure to entice players to try this game."

    # Now remove the bold "Play Aztec Bonanza Free..." paragraph entirely.
    $boldPara.Range.Delete()

    Write-Output "done"
}
